# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja "tasas": actualizar tasas de cambio ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 68.45
$wsTasas.Range("O10").Value = 4110.63
$wsTasas.Range("N12").Value = 4137.99
$wsTasas.Range("O12").Value = 67.39

# --- Hoja "Hoja1": actualizar texto de conversión del día ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 14.61 = 60053.03 pesos`n✅ 60053.03 pesos = 14.51 = 978.0 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
